# "Run first successful test."
# All boolean test-result cells on the "Test Results" sheet flip from
# FALSE to TRUE, reflecting the first fully-passing test run. The user
# had also scrolled / clicked down into the grid (around row 13) when
# they made the edit, so the sheet's active selection is updated too.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Results")
$ws.Activate()

# Rows 3-15 (Z02-Z09, C01-C05): only the "Delete Test Passed" column (E)
# was still FALSE - flip it to TRUE now that it passes too.
$ws.Range("E3:E15").Value = $true

# Rows 16-24 (D01-D09): none of the Create/Read/Update/Delete tests had
# been run yet (all FALSE) - this is the first successful run for the
# Device rows, so every column passes now.
$ws.Range("B16:E24").Value = $true

# Reflect where the user ended up after scrolling/clicking through the
# freshly-updated rows.
[void]$ws.Range("H13").Select()
